# TestBrugerJonstrup.xlsx — "PSADT and Packages examples"
#
# The sheet held a demo list of users including a "Mobil nummer" (mobile
# phone number) column (F). This edit scrubs those sample phone numbers
# out of the sheet (keeping the header and every other column intact),
# and leaves the selection parked on the now-empty F2:F9 range, matching
# where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Clear out the sample mobile numbers in F2:F9. ClearContents() removes the
# stored value but (for the cells that already carried explicit formatting,
# F3/F5:F9) keeps the cell's style around as an empty formatted cell - same
# as F2/F4 which had no extra formatting and end up fully empty again.
$ws.Range("F2:F9").ClearContents()

# Column D (MA-nummer) used a style that was identical in every respect to
# the one already used elsewhere on the sheet (e.g. B2) except for a
# redundant "applyNumberFormat" flag. Normalize D2:D9 onto that shared
# vertical-center style so the workbook doesn't carry a duplicate.
$ws.Range("D2:D9").VerticalAlignment = -4108  ## xlCenter

# The built-in hyperlink cell style shipped under its Danish/legacy label
# ("Link"); rename it to the standard "Hyperlink" name.
try {
    $linkStyle = $wb.Styles.Item("Link")
    $linkStyle.Name = "Hyperlink"
} catch {
    # Older/limited hosts may not expose named-style renaming; non-fatal.
}

# Leave the selection on the column that was just edited.
$ws.Range("F2:F9").Select() | Out-Null
